# Insurance sheet ("保險" / sheet4): normalize rows to add the common
# metadata columns (category, date, legislator_name, legislator_id,
# source_file, index) and rewrite the B/C/E columns to use a shared
# "company" / policy-name scheme, collapsing the long per-row insurance
# period/premium strings down to the constant "insurance" category value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# "date" column (G) holds the literal text "2011-11-21" in every other
# sheet of this workbook (shared string, not a real Excel date) - force
# text formatting first so COM doesn't silently coerce it to a date serial.
$ws.Range("G1:G9").NumberFormat = "@"

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# ---- Data rows (2-9) ----
# row, A(index no.), B(company), C(policy name), D(owner)
$rows = @(
    @(2, 87, "南山人壽", "子女教育保險",             "廖述嘉"),
    @(3, 88, "中華郵政", "十年快樂兒同還本終身壽險", "盧秀燕"),
    @(4, 89, "中華郵政", "十年快樂兒同還本終身壽險", "盧秀燕"),
    @(5, 90, "中華郵政", "金寶貝兒童保險",           "盧秀燕"),
    @(6, 91, "中國人壽", "得意人生終身保險",         "盧秀燕"),
    @(7, 92, "中國人壽", "得意人生終身保險",         "盧秀燕"),
    @(8, 93, "中國人壽", "得意人生終身保險",         "盧秀燕"),
    @(9, 94, "中國人壽", "得意人生終身保險",         "盧秀燕")
)

foreach ($row in $rows) {
    $r       = $row[0]
    $idx     = $row[1]
    $company = $row[2]
    $name    = $row[3]
    $owner   = $row[4]

    $ws.Range("B$r").Value = $company
    $ws.Range("C$r").Value = $name
    $ws.Range("D$r").Value = $owner
    $ws.Range("E$r").Value = "insurance"

    $ws.Range("F$r").Value = "normal"
    $ws.Range("G$r").Value = "2011-11-21"
    $ws.Range("H$r").Value = "盧秀燕"
    $ws.Range("I$r").Value = 869
    $ws.Range("J$r").Value = "tmp9eb41"
    $ws.Range("K$r").Value = $idx
}
